$d = $word.ActiveDocument

# 1) Apply the new paragraph indent (left=180 twips=9pt, right=206 twips=10.3pt)
#    to every paragraph in the document (the 9 existing NormalWeb paragraphs
#    plus the trailing empty paragraph that only carried the _GoBack bookmark).
$d.Paragraphs.LeftIndent = 9
$d.Paragraphs.RightIndent = 10.3

# 2) Move the "_GoBack" bookmark from the trailing empty paragraph to sit
#    right after the last run of the "Findings" paragraph (immediately after
#    "...abdominal wall."), collapsed (zero-length), matching the target XML.
#    A directly-collapsed Range confuses Bookmarks.Add in this host, so we
#    insert a one-character placeholder, bookmark that, then delete the
#    placeholder text while the bookmark marks stay put.
$findRng = $d.Content
$findRng.Find.Execute("abdominal wall.", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$findRng.Collapse(0)
$findRng.InsertAfter("X")
$d.Bookmarks.Add("_GoBack", $findRng)
$placeholder = $d.Range($findRng.Start, $findRng.Start + 1)
$placeholder.Text = ""

# 3) Update the section page margins: right 1800 -> 746 twips (37.3pt),
#    left 1800 -> 810 twips (40.5pt).
$d.PageSetup.RightMargin = 37.3
$d.PageSetup.LeftMargin = 40.5

Write-Output "edits applied"
